$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.002.96"

$ws.Range("D3").Value = "1.556.11"
$ws.Range("E3").Value = "  +0.75%  "

$ws.Range("E4").Value = "  -0.22%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.43%  "

$ws.Range("E6").Value = "  +0.21%  "

$ws.Range("E7").Value = "  -0.26%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.06"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.78%  "

$ws.Range("E9").Value = "  +0.11%  "

$ws.Range("E10").Value = "  +0.99%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0857"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.09%  "

$ws.Range("D12").Value = "1.777.34"
$ws.Range("E12").Value = "  +0.72%  "

$ws.Range("D13").Value = "1.555.80"
$ws.Range("E13").Value = "  +0.50%  "

$ws.Range("E14").Value = "  +1.44%  "

$ws.Range("E15").Value = "  +1.72%  "

$ws.Range("D16").Value = "26.994.21"
$ws.Range("E16").Value = "  +0.44%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.81"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.74%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "218.22"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.04%  "

$ws.Range("D19").Value = "0.0₃0695"
$ws.Range("E19").Value = "  +1.84%  "

$ws.Range("E20").Value = "  +1.99%  "

$ws.Range("E21").Value = "  -0.25%  "

$ws.Range("E22").Value = "  +1.43%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.24"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.94%  "

$ws.Range("E24").Value = "  +0.78%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.34"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.38%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.67"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.14%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.01%  "

$ws.Range("E28").Value = "  +1.05%  "

$ws.Range("E29").Value = "  -0.22%  "

$ws.Range("E30").Value = "  +2.19%  "

$ws.Range("E31").Value = "  -0.26%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.24"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.75%  "

$ws.Range("D33").Value = "1.424.02"
$ws.Range("E33").Value = "  +4.90%  "

$ws.Range("E34").Value = "  +4.77%  "

$ws.Range("E35").Value = "  +4.07%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.982"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.45%  "

$ws.Range("E37").Value = "  +0.17%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0166"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.93%  "

$ws.Range("E39").Value = "  +0.76%  "

$ws.Range("E40").Value = "  +1.27%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.71"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.64%  "

$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.21%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.32"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.93%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.985"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.49%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.58"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.06%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.76"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.68%  "

$ws.Range("D47").Value = "1.691.16"
$ws.Range("E47").Value = "  +0.68%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.99"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.75%  "

$ws.Range("E49").Value = "  +1.91%  "

$ws.Range("D50").Value = "0.0₇0997"
$ws.Range("E50").Value = "  +1.96%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0956"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.11%  "
